# issue #5: add legislator_id, name, date into dataframe
#
# The "股票" (stocks) worksheet (4th sheet) gets three new trailing
# columns - date / legislator_name / legislator_id - appended after the
# existing "total" column, for every existing data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

$legislatorName = "葉津鈴"
$legislatorId = 1805
$reportDate = "2013-09-15"

# --- Header row (row 1): copy the formatting of the last existing header
# cell (G1, style "1" - bold + border) onto the three new header cells,
# then stamp in the column names.
$ws.Range("G1").Copy()
$ws.Range("H1:J1").PasteSpecial(-4122)
$ws.Range("H1").Value = "date"
$ws.Range("I1").Value = "legislator_name"
$ws.Range("J1").Value = "legislator_id"

# --- Data rows: for every row already populated in column A (the record
# id column), copy the formatting of the existing "total" cell (G, style
# "2") across into H:J, then fill in the date/name/id values. The date
# column is forced to Text format first so "2013-09-15" is stored as a
# literal string rather than being reinterpreted as a date serial.
$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Range("G$r").Copy()
    $ws.Range("H$r`:J$r").PasteSpecial(-4122)

    $ws.Range("H$r").NumberFormat = "@"
    $ws.Range("H$r").Value = $reportDate
    $ws.Range("I$r").Value = $legislatorName
    $ws.Range("J$r").Value = $legislatorId
}

$excel.CutCopyMode = $false
